$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume snapshot data refreshed by GitHub Actions
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.27"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.12"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.399"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05996"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.428"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.508"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8110"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9231"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1424"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07404"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03259"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03072"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09342"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.846"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001582"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04701"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005890"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005863"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004881"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006798"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.566"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.137"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002340"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03968"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006598"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004399"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1077"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009193"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005091"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7000"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
